# Refactor synthetic array: insert a new "statut_name" column (French status
# label) right after "statut_label", shifting NCTId..intervention_type one
# column to the right (C -> D ... L -> M).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at C; everything from C onward (including all
# formatting) shifts right by one column.
$ws.Columns("C").Insert()

# New header for the inserted column.
$ws.Range("C1").Value = "statut_name"

# Per-row French status-name values for the new column.
$ws.Range("C2").Value = "résultat et / ou publication posté"
$ws.Range("C3").Value = "résultat et / ou publication posté"
$ws.Range("C4").Value = "pas de résultat ni de publication"
$ws.Range("C5").Value = "résultat et / ou publication posté"
$ws.Range("C6").Value = "résultat et / ou publication posté"
$ws.Range("C7").Value = "résultat et / ou publication posté dans les 36 mois"
$ws.Range("C8").Value = "résultat et / ou publication posté dans les 12 mois"
$ws.Range("C9").Value = "résultat et / ou publication posté dans les 12 mois"
